$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need to be forced
# to text so Excel does not auto-convert the inline string into a Number
# (these columns hold formatted price/volume text, not numeric values).
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '28.494.66'
$ws.Range("E2").Value = '  +4.31%  '
$ws.Range("D3").Value = '1.592.32'
$ws.Range("E4").Value = '  -0.33%  '
Set-TextValue "D5" '214.93'
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("E7").Value = '  -0.28%  '
Set-TextValue "D8" '24.01'
$ws.Range("E8").Value = '  +8.16%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("D12").Value = '1.820.18'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").Value = '1.589.71'
$ws.Range("E13").Value = '  +1.42%  '
Set-TextValue "D14" '3.80'
$ws.Range("E14").Value = '  +0.42%  '
Set-TextValue "D15" '0.532'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '28.489.89'
$ws.Range("E16").Value = '  +4.52%  '
$ws.Range("E17").Value = '  +2.78%  '
Set-TextValue "D18" '232.59'
$ws.Range("E18").Value = '  +7.03%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0710'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D20" '7.51'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("E24").Value = '  +1.38%  '
Set-TextValue "D25" '151.95'
$ws.Range("E25").Value = '  -1.16%  '
Set-TextValue "D26" '15.30'
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").Value = '1.413.88'
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("E36").Value = '  -5.22%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  +0.28%  '
Set-TextValue "D39" '2.55'
$ws.Range("E39").Value = '  +9.00%  '
Set-TextValue "D40" '0.543'
$ws.Range("E40").Value = '  +1.52%  '
$ws.Range("E41").Value = '  +0.54%  '
Set-TextValue "D42" '5.78'
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("E45").Value = '  +4.83%  '
Set-TextValue "D46" '64.74'
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = '1.733.06'
$ws.Range("E47").Value = '  +1.63%  '
Set-TextValue "D48" '87.97'
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("E49").Value = '  +4.65%  '
$ws.Range("E50").Value = '  -0.42%  '
Set-TextValue "D51" '39.64'
$ws.Range("E51").Value = '  +16.74%  '
